$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date headers (H1, J1, L1) with their merged partner cells (I1, K1, M1) ---
# Copy formatting from the existing F1:G1 pair (date header + merged blank) onto
# the three new pairs, then set the new dates.
$ws.Range("F1:G1").Copy() | Out-Null
$ws.Range("H1:I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F1:G1").Copy() | Out-Null
$ws.Range("J1:K1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1:G1").Copy() | Out-Null
$ws.Range("L1:M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("H1").Value = 41730
$ws.Range("J1").Value = 41731
$ws.Range("L1").Value = 41732

# --- Merge the new header pairs ---
$ws.Range("H1:I1").Merge() | Out-Null
$ws.Range("J1:K1").Merge() | Out-Null
$ws.Range("L1:M1").Merge() | Out-Null

# --- New data values for the extra three days ---
$ws.Range("H3").Value = 20
$ws.Range("J3").Value = 21
$ws.Range("L3").Value = 21

$ws.Range("H4").Value = 13
$ws.Range("J4").Value = 13
$ws.Range("L4").Value = 13

$ws.Range("H5").Value = 10
$ws.Range("J5").Value = 11
$ws.Range("L5").Value = 11

$ws.Range("H6").Value = 22
$ws.Range("J6").Value = 22
$ws.Range("L6").Value = 22

$ws.Range("D7").Value = 11
$ws.Range("H7").Value = 11
$ws.Range("J7").Value = 11
$ws.Range("L7").Value = 11

# --- New column widths ---
$ws.Columns.Item(9).ColumnWidth = 18.140625
$ws.Columns.Item(10).ColumnWidth = 9.42578125
$ws.Columns.Item(11).ColumnWidth = 18.5703125
$ws.Columns.Item(13).ColumnWidth = 14.7109375

# --- View changes ---
$window = $excel.ActiveWindow
$window.Zoom = 85
$window.ScrollColumn = 3
$ws.Range("L8").Select() | Out-Null
